$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range('D2').Value = '27.047.00'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '1.826.36'
$ws.Range('E3').Value = '  +0.03%  '
Set-TextValue $ws 'D4' '1.0000'
$ws.Range('E4').Value = '  -0.33%  '
Set-TextValue $ws 'D5' '311.83'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('E6').Value = '  -0.25%  '
Set-TextValue $ws 'D7' '0.4409'
$ws.Range('E7').Value = '  +2.67%  '
Set-TextValue $ws 'D8' '0.3683'
$ws.Range('E8').Value = '  -0.39%  '
Set-TextValue $ws 'D9' '0.07272'
$ws.Range('E9').Value = '  +0.38%  '
Set-TextValue $ws 'D10' '0.8447'
$ws.Range('E10').Value = '  -2.29%  '
Set-TextValue $ws 'D11' '20.72'
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('D12').Value = '1.815.28'
$ws.Range('E12').Value = '  -0.58%  '
Set-TextValue $ws 'D13' '6.661'
$ws.Range('E13').Value = '  -0.07%  '
Set-TextValue $ws 'D14' '0.07069'
$ws.Range('E14').Value = '  -0.27%  '
Set-TextValue $ws 'D15' '5.303'
$ws.Range('E15').Value = '  -0.78%  '
Set-TextValue $ws 'D16' '89.83'
$ws.Range('E16').Value = '  +2.51%  '
Set-TextValue $ws 'D17' '1.001'
$ws.Range('E17').Value = '  -0.31%  '
Set-TextValue $ws 'D18' '0.000008798'
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range('E19').Value = '  -0.20%  '
Set-TextValue $ws 'D20' '14.93'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').Value = '27.029.11'
$ws.Range('E21').Value = '  -1.21%  '
Set-TextValue $ws 'D22' '5.155'
$ws.Range('E22').Value = '  -0.04%  '
Set-TextValue $ws 'D23' '10.90'
$ws.Range('E23').Value = '  +0.46%  '
$ws.Range('D24').Value = '2.046.24'
$ws.Range('E24').Value = '  -0.25%  '
Set-TextValue $ws 'D25' '1.984'
$ws.Range('E25').Value = '  -1.24%  '
Set-TextValue $ws 'D26' '151.64'
$ws.Range('E26').Value = '  -0.97%  '
Set-TextValue $ws 'D27' '2.209'
$ws.Range('E27').Value = '  +3.27%  '
Set-TextValue $ws 'D28' '18.29'
$ws.Range('E28').Value = '  -0.87%  '
Set-TextValue $ws 'D29' '5.237'
$ws.Range('E29').Value = '  -0.94%  '
Set-TextValue $ws 'D30' '117.04'
$ws.Range('E30').Value = '  +0.29%  '
Set-TextValue $ws 'D31' '0.08786'
$ws.Range('E31').Value = '  -0.55%  '
Set-TextValue $ws 'D32' '1.179'
$ws.Range('E32').Value = '  -1.92%  '
Set-TextValue $ws 'D33' '0.7415'
$ws.Range('E33').Value = '  -3.06%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D34' '4.427'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws 'D35' '2.889'
$ws.Range('E35').Value = '  +1.36%  '
Set-TextValue $ws 'D36' '0.9997'
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('E37').Value = '  -2.31%  '
Set-TextValue $ws 'D38' '0.01948'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('E39').Value = '  -0.40%  '
Set-TextValue $ws 'D40' '7.258'
$ws.Range('E40').Value = '  +2.12%  '
$ws.Range('E41').Value = '  -0.48%  '
Set-TextValue $ws 'D42' '0.5167'
$ws.Range('E42').Value = '  +2.19%  '
Set-TextValue $ws 'D43' '0.1698'
$ws.Range('E43').Value = '  +1.17%  '
Set-TextValue $ws 'D44' '8.541'
$ws.Range('E44').Value = '  -1.19%  '
Set-TextValue $ws 'D45' '10.62'
$ws.Range('E45').Value = '  +0.58%  '
Set-TextValue $ws 'D46' '0.4818'
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('E47').Value = '  -0.14%  '
Set-TextValue $ws 'D48' '1.933'
$ws.Range('E48').Value = '  +6.38%  '
Set-TextValue $ws 'D49' '0.9994'
$ws.Range('E49').Value = '  -0.31%  '
$ws.Range('E50').Value = '  -1.40%  '
Set-TextValue $ws 'D51' '1.660'
$ws.Range('E51').Value = '  -0.37%  '
